$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: turn the "time of day" decimal values (custom time format) into
# explicit "HH:MM~HH:MM" text range labels ------------------------------------
$ws.Range("A4").Value = "13:00~14:00"
$ws.Range("A5").Value = "14:00~15:00"
$ws.Range("A6").Value = "15:00~16:00"
$ws.Range("A7").Value = "16:00~17:00"
$ws.Range("A8").Value = "17:00~18:00"
$ws.Range("A9").Value = "18:00~19:00"
$ws.Range("A10").Value = "19:00~20:00"

# --- Add the "kiwoom"/영어 3rd day: Tuesday (D) and Thursday' (H) now mirror
# what Monday (B) / Friday (I) already have for the 16:00 and 18:00 slots -----

# Row 7 (16:00~17:00): 영어(16시 40분) on Tuesday and Thursday'
$ws.Range("B7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$ws.Range("D7").Value = $ws.Range("B7").Value()

$ws.Range("B7").Copy() | Out-Null
$ws.Range("H7").PasteSpecial(-4122) | Out-Null
$ws.Range("H7").Value = $ws.Range("B7").Value()

# Row 8 (17:00~18:00): keep the same highlighted (empty) block going
$ws.Range("B8").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null

$ws.Range("B8").Copy() | Out-Null
$ws.Range("H8").PasteSpecial(-4122) | Out-Null

# Row 9 (18:00~19:00): 18시 40분 on Tuesday and Thursday'
$ws.Range("B9").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Value = $ws.Range("B9").Value()

$ws.Range("B9").Copy() | Out-Null
$ws.Range("H9").PasteSpecial(-4122) | Out-Null
$ws.Range("H9").Value = $ws.Range("B9").Value()

$excel.CutCopyMode = $false

# --- Column widths: column A becomes visible with its own width, column B
# (still hidden) grows slightly ------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10.14285714285714
$ws.Columns.Item(2).ColumnWidth = 14.71428571428571

# --- Selection bookmark moves from K15 to F15 --------------------------------
$ws.Range("F15").Select() | Out-Null
